$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (25) with the latest Adafruit IO reading, mirroring
# the data already present in row 24.
$row = 25

$ws.Range("A$row").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$row").Value = "temperature"

# Column C holds a numeric-looking reading ("25") that must stay stored
# as text, just like the rest of the sheet (all inline/shared strings).
# Temporarily force a text format so Excel doesn't coerce it to a
# number, then clear the formatting again so no stray style is left
# behind on the cell.
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = "25"
$ws.Range("C$row").ClearFormats()

$ws.Range("D$row").Value = "N/A"
$ws.Range("E$row").Value = "N/A"
$ws.Range("F$row").Value = "N/A"
